# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a handful of rows whose tags changed
# after the transcripts were cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 10; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 17; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 21; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 31; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 34; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
